$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(0.129, -0.033, 0.0667, 0.1578645893000344, 0.1578645893000344, 0.07520408622817321, 0.06252432750365564, 2288, 0.06552869744529728, 2035, 0.02929452732164527, 0.8894230769230769, 1375, 0.01979359954165221, 0.6009615384615384, 660, 0.3243243243243243, 1707, 0.02457285412189115, 0.04192626255222458, 0.04789904570080378, -0.005972783148579193, 0.6411208561022858, 0.04008565037636339, 0.04112498518784301, -0.001039334811479622, 19702, 615.8706262855198, 20317.87062628552, 18610.87062628552, 0.2262952891070507, 0.2647939539907328, 0.2113004279507885, 0.2480659884770725, 524, 524, 6.476660092044708, 4.944656488549619, 6.117971935005102, 4.944656488549619)

for ($i = 0; $i -lt $values.Length; $i++) {
    $col = 4 + $i  # column D = 4
    $ws.Cells.Item(2, $col).Value = $values[$i]
    $ws.Cells.Item(3, $col).Value = $values[$i]
}
